# excel_to_ruby now does multiple passes, converting to values then
# replacing indirects -- add a worked example on the "Referencing" sheet
# showing INDIRECT() built dynamically from two cell values.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Referencing")

$ws.Range("B11").Value = "Named"
$ws.Range("C11").Value = "Reference"
$ws.Range("D11").Formula = '=INDIRECT(B11&"_"&C11)'

# Leave the selection where Excel would land after typing the formula
# into D11 and hitting Enter.
$ws.Range("D12").Select()
